{"js": "// Planning.docx content edit:\n//  1. \"Female\" -> \"Women\" and \"Male\" -> \"Men\" everywhere (whole word, case-sensitive)\n//     this covers: the \"Gender type\" options paragraph (2 occurrences each),\n//     the \"Female/Male out of order\" paragraph labels, and the\n//     \"Female/Male Updated: date\" paragraph labels.\n//  2. Remove the duplicate \"HasUrinals: dropdown, options: Yes, No, Unknown, required\"\n//     paragraph (a leftover duplicate of the earlier HasUrinals field) - this was the\n//     \"fix for variable creation\" mentioned in the commit message; removing it shifts\n//     HasShowers/ADA_Accessible/Usage_Fee/DamageDescription/MessDescription/\n//     ShortageDescription back into their intended single slots.\n\n// Step 1: word replacements.\nconst femaleResults = context.document.body.search(\"Female\", {\n  matchCase: true,\n  matchWholeWord: true,\n});\nfemaleResults.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < femaleResults.items.length; i++) {\n  femaleResults.items[i].insertText(\"Women\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\nconst maleResults = context.document.body.search(\"Male\", {\n  matchCase: true,\n  matchWholeWord: true,\n});\nmaleResults.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < maleResults.items.length; i++) {\n  maleResults.items[i].insertText(\"Men\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// Step 2: delete the duplicate \"HasUrinals\" paragraph. Find it by searching for the\n// exact paragraph text (there are two matching paragraphs before the fix; the first\n// one is the legitimate field, the second is the duplicate that must be removed).\nconst dupResults = context.document.body.search(\n  \"HasUrinals: dropdown, options: Yes, No, Unknown, required\",\n  { matchCase: true }\n);\ndupResults.load(\"items\");\nawait context.sync();\n\nif (dupResults.items.length > 1) {\n  const dupParagraph = dupResults.items[1].paragraphs.getFirst();\n  dupParagraph.delete();\n  await context.sync();\n}\n", "ps1": "# Planning.docx content edit:\n#  1. \"Female\" -> \"Women\" and \"Male\" -> \"Men\" everywhere (whole word, case-sensitive).\n#     This covers: the \"Gender type\" options paragraph (2 occurrences each),\n#     the \"Female/Male out of order\" paragraph labels, and the\n#     \"Female/Male Updated: date\" paragraph labels.\n#  2. Remove the duplicate \"HasUrinals: dropdown, options: Yes, No, Unknown, required\"\n#     paragraph (a leftover duplicate of the earlier HasUrinals field) - this was the\n#     \"fix for variable creation\" mentioned in the commit message; removing it shifts\n#     HasShowers/ADA_Accessible/Usage_Fee/DamageDescription/MessDescription/\n#     ShortageDescription back into their intended single slots.\n\n$d = $word.ActiveDocument\n\n# Step 1: word replacements.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"Female\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $true\n$find.Replacement.Text = \"Women\"\n$find.Execute($null, $true, $true, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Text = \"Male\"\n$find2.MatchCase = $true\n$find2.MatchWholeWord = $true\n$find2.Replacement.Text = \"Men\"\n$find2.Execute($null, $true, $true, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n\n# Step 2: locate the second occurrence of the duplicated \"HasUrinals\" paragraph text\n# and delete that whole paragraph.\n$searchRange = $d.Content\n$searchRange.Find.ClearFormatting()\n$searchRange.Find.Text = \"HasUrinals: dropdown, options: Yes, No, Unknown, required\"\n$searchRange.Find.MatchCase = $true\n$occurrence = 0\n$dupStart = -1\nwhile ($searchRange.Find.Execute()) {\n    $occurrence++\n    if ($occurrence -eq 2) {\n        $dupStart = $searchRange.Start\n        break\n    }\n    $searchRange.Collapse(0)\n}\n\nif ($dupStart -ge 0) {\n    $targetIndex = -1\n    $i = 1\n    foreach ($p in $d.Paragraphs) {\n        if ($p.Range.Start -eq $dupStart) {\n            $targetIndex = $i\n            break\n        }\n        $i++\n    }\n    if ($targetIndex -ge 1) {\n        $d.Paragraphs.Item($targetIndex).Range.Delete()\n    }\n}\n"}
